# Add instruction to write -1 if you give up on a challenge.
# On the "Battle" sheet, insert two new rows above row 12 (shifting the
# existing question/merged-range block down by two rows) and write the new
# instruction text into the newly created B12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Battle")

$ws.Rows("12:13").Insert()

$ws.Range("B12").Value = "If you give up, write '-1' as the answer."
